# REF/SM: Generate settings for LS3 cavity tuning.
#
# The sheet gained a new (empty) column B: every existing column from B
# onward (the "[MeV/u]" header plus all the magnet-setting columns, for
# both the "LS3 transport" and "LS3 acceleration" rows) shifted one
# column to the right, column A ("LS3 transport"/"LS3 acceleration"
# labels) stayed put, and the blank column B was left with default
# width/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column B - shifts B:BR -> C:BS.
$ws.Columns("B:B").Insert()

# Matches the saved view state (active cell moved to C2).
$ws.Range("C2").Select()
